$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 620.8182
$ws.Range("I2").Value = 620.8182
$ws.Range("K2").Value = 620.8182
$ws.Range("M2").Value = -507.8182

$ws.Range("H32").Value = 932.5
$ws.Range("I32").Value = 400
$ws.Range("J32").Value = 980.9091
$ws.Range("K32").Value = 400
$ws.Range("L32").Value = 980.9091
$ws.Range("M32").Value = -74
$ws.Range("N32").Value = -1632.9091

$ws.Range("H99").Value = 464.7143
$ws.Range("J99").Value = 241
$ws.Range("L99").Value = 723
$ws.Range("N99").Value = -3719

$ws.Range("H107").Value = 7578100
$ws.Range("I107").Value = 8929379
$ws.Range("J107").Value = 10936
$ws.Range("K107").Value = 8929379
$ws.Range("L107").Value = 10936
$ws.Range("M107").Value = -8927459
$ws.Range("N107").Value = -14776

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1564.0526
$ws.Range("I2").Value = 1732.125
$ws.Range("J2").Value = 667.6667
$ws.Range("K2").Value = 1732.125
$ws.Range("L2").Value = 667.6667
$ws.Range("M2").Value = -1619.125
$ws.Range("N2").Value = -893.6667

$ws.Range("H45").Value = 10544.4375
$ws.Range("I45").Value = 11771.214
$ws.Range("J45").Value = 1957
$ws.Range("K45").Value = 11771.214
$ws.Range("L45").Value = 1957
$ws.Range("M45").Value = -11394.214
$ws.Range("N45").Value = -2711

$ws.Range("H61").Value = 3962.5
$ws.Range("I61").Value = 5749.926
$ws.Range("K61").Value = 5749.926
$ws.Range("M61").Value = -5537.926

$ws.Range("H116").Value = 1564.0526
$ws.Range("I116").Value = 1732.125
$ws.Range("J116").Value = 667.6667
$ws.Range("K116").Value = 1732.125
$ws.Range("L116").Value = 667.6667
$ws.Range("M116").Value = 561.875
$ws.Range("N116").Value = -5255.6667

$ws.Range("H136").Value = 3962.5
$ws.Range("I136").Value = 5749.926
$ws.Range("K136").Value = 17249.778
$ws.Range("M136").Value = -14699.778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1564.0526
$ws.Range("I3").Value = 1732.125
$ws.Range("J3").Value = 667.6667
$ws.Range("K3").Value = 1732.125
$ws.Range("L3").Value = 667.6667
$ws.Range("M3").Value = -1618.125
$ws.Range("N3").Value = -895.6667

$ws.Range("H94").Value = 1575.4667
$ws.Range("I94").Value = 1269.8096
$ws.Range("J94").Value = 2288.6667
$ws.Range("K94").Value = 1269.8096
$ws.Range("L94").Value = 2288.6667
$ws.Range("M94").Value = -818.8096
$ws.Range("N94").Value = -3190.6667

$ws.Range("H99").Value = 76925440
$ws.Range("I99").Value = 142859710
$ws.Range("J99").Value = 2116.6667
$ws.Range("K99").Value = 142859710
$ws.Range("L99").Value = 2116.6667
$ws.Range("M99").Value = -142858212
$ws.Range("N99").Value = -5112.6667

$ws.Range("H132").Value = 35000
$ws.Range("J132").Value = 35000
$ws.Range("L132").Value = 35000
$ws.Range("N132").Value = -45120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4343.409
$ws.Range("I62").Value = 3363.8333
$ws.Range("J62").Value = 8751.5
$ws.Range("K62").Value = 3363.8333
$ws.Range("L62").Value = 8751.5
$ws.Range("M62").Value = -2739.8333
$ws.Range("N62").Value = -9999.5

$ws.Range("H65").Value = 4343.409
$ws.Range("I65").Value = 3363.8333
$ws.Range("J65").Value = 8751.5
$ws.Range("K65").Value = 16819.1665
$ws.Range("L65").Value = 43757.5
$ws.Range("M65").Value = -13699.1665
$ws.Range("N65").Value = -49997.5

$ws.Range("H122").Value = 1095.381
$ws.Range("I122").Value = 1153.4615
$ws.Range("J122").Value = 1001
$ws.Range("K122").Value = 3460.3845
$ws.Range("L122").Value = 3003
$ws.Range("M122").Value = -1010.3845
$ws.Range("N122").Value = -7903

$ws.Range("H132").Value = 3521.7778
$ws.Range("I132").Value = 2642
$ws.Range("K132").Value = 7926
$ws.Range("M132").Value = -5396

$ws.Range("H134").Value = 2602.25
$ws.Range("I134").Value = 2841.3572
$ws.Range("K134").Value = 8524.071599999999
$ws.Range("M134").Value = -5989.071599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 321.91306
$ws.Range("I40").Value = 158
$ws.Range("J40").Value = 1100.5
$ws.Range("K40").Value = 632
$ws.Range("L40").Value = 4402
$ws.Range("M40").Value = -563
$ws.Range("N40").Value = -4540

$ws.Range("H68").Value = 2528.8489
$ws.Range("I68").Value = 2883.2979
$ws.Range("J68").Value = 2101.6924
$ws.Range("K68").Value = 8649.893700000001
$ws.Range("L68").Value = 6305.0772
$ws.Range("M68").Value = -7838.893700000001
$ws.Range("N68").Value = -7927.0772

$ws.Range("H71").Value = 2528.8489
$ws.Range("I71").Value = 2883.2979
$ws.Range("J71").Value = 2101.6924
$ws.Range("K71").Value = 25949.6811
$ws.Range("L71").Value = 18915.2316
$ws.Range("M71").Value = -21893.6811
$ws.Range("N71").Value = -27027.2316

$ws.Range("H129").Value = 2081
$ws.Range("J129").Value = 2826.0715
$ws.Range("L129").Value = 8478.2145
$ws.Range("N129").Value = -18478.2145

$ws.Range("H133").Value = 33643.36
$ws.Range("I133").Value = 115787.89
$ws.Range("K133").Value = 347363.67
$ws.Range("M133").Value = -342303.67

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2862.7144
$ws.Range("I102").Value = 1006
$ws.Range("J102").Value = 4255.25
$ws.Range("K102").Value = 1006
$ws.Range("L102").Value = 4255.25
$ws.Range("M102").Value = 616
$ws.Range("N102").Value = -7499.25

$ws.Range("H113").Value = 58824500
$ws.Range("I113").Value = 71429540
$ws.Range("J113").Value = 996.6667
$ws.Range("K113").Value = 71429540
$ws.Range("L113").Value = 996.6667
$ws.Range("M113").Value = -71427370
$ws.Range("N113").Value = -5336.6667

$ws.Range("H126").Value = 8467.333000000001
$ws.Range("I126").Value = 9446.923000000001
$ws.Range("J126").Value = 2100
$ws.Range("K126").Value = 28340.769
$ws.Range("L126").Value = 6300
$ws.Range("M126").Value = -25870.769
$ws.Range("N126").Value = -11240

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 567.1
$ws.Range("J16").Value = 625
$ws.Range("L16").Value = 625
$ws.Range("N16").Value = -965
